$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33 / 34 swap: Aptos <-> FirstDigitalUSD ---
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'0.997"
$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'5.70"
$ws.Range("E34").Value = "  -1.32%  "

# --- Price / Volume updates for remaining rows ---
$ws.Range("D2").Value = "53.604.12"
$ws.Range("E2").Value = "  -5.09%  "
$ws.Range("D3").Value = "2.225.08"
$ws.Range("E3").Value = "  -6.51%  "
$ws.Range("D5").Value = "'483.96"
$ws.Range("E5").Value = "  -3.86%  "
$ws.Range("D6").Value = "'125.85"
$ws.Range("E6").Value = "  -3.63%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'0.517"
$ws.Range("E8").Value = "  -5.37%  "
$ws.Range("D9").Value = "2.247.32"
$ws.Range("E9").Value = "  -5.95%  "
$ws.Range("D10").Value = "'0.0916"
$ws.Range("E10").Value = "  -7.22%  "
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("D12").Value = "'4.70"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("E13").Value = "  -3.71%  "
$ws.Range("D14").Value = "2.623.99"
$ws.Range("E14").Value = "  -6.44%  "
$ws.Range("D15").Value = "'21.02"
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("D16").Value = "53.495.08"
$ws.Range("E16").Value = "  -5.26%  "
$ws.Range("E17").Value = "  -3.95%  "
$ws.Range("D18").Value = "2.236.76"
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("D19").Value = "'9.55"
$ws.Range("E19").Value = "  -5.06%  "
$ws.Range("D20").Value = "'3.94"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Value = "'298.05"
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("D22").Value = "'6.09"
$ws.Range("E22").Value = "  -2.73%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'63.38"
$ws.Range("E24").Value = "  -3.32%  "
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "'0.362"
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("D27").Value = "'0.142"
$ws.Range("E27").Value = "  -4.43%  "
$ws.Range("D28").Value = "'6.98"
$ws.Range("E28").Value = "  -4.84%  "
$ws.Range("D29").Value = "'169.64"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").Value = "0.0₃0681"
$ws.Range("E30").Value = "  -5.13%  "
$ws.Range("D31").Value = "'1.58"
$ws.Range("E31").Value = "  -3.61%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E35").Value = "  -3.88%  "
$ws.Range("D36").Value = "'17.42"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("E38").Value = "  +5.11%  "
$ws.Range("E39").Value = "  -6.21%  "
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").Value = "'0.364"
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("E42").Value = "  -3.06%  "
$ws.Range("D43").Value = "'3.27"
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("D44").Value = "'122.23"
$ws.Range("E44").Value = "  -6.79%  "
$ws.Range("D45").Value = "'4.63"
$ws.Range("E45").Value = "  -7.05%  "
$ws.Range("D46").Value = "'0.0876"
$ws.Range("E46").Value = "  -3.67%  "
$ws.Range("D47").Value = "'0.535"
$ws.Range("E47").Value = "  -5.43%  "
$ws.Range("D48").Value = "'230.95"
$ws.Range("E48").Value = "  -4.61%  "
$ws.Range("D49").Value = "'0.0469"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("D50").Value = "'0.0201"
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("D51").Value = "'15.97"
$ws.Range("E51").Value = "  -5.63%  "
